$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project details")

$values = @(
    "Solar Project DA",
    "Solar Project DB",
    "Solar Project DA",
    "Solar Project AA",
    "Solar Project AD",
    "Solar Project DA",
    "Solar Project AB",
    "Solar Project BA",
    "Solar Project CD",
    "Solar Project DC",
    "Solar Project DA",
    "Solar Project DB",
    "Solar Project BA",
    "Solar Project AA",
    "Solar Project CB",
    "Solar Project CC",
    "Solar Project BB",
    "Solar Project BB",
    "Solar Project BC",
    "Solar Project BA",
    "Solar Project BA",
    "Solar Project DB",
    "Solar Project AC",
    "Solar Project CB"
)

$startRow = 40
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
